# Update the "想去人数" (F column) counts on the "展览" and "全部类型" sheets
# to reflect the latest generated output.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 3786
    4  = 2278
    5  = 445
    6  = 8
    7  = 17
    9  = 105
    10 = 88
    11 = 1403
    12 = 247
    13 = 2330
    14 = 165
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $updates[$row]
}

$updates4 = @{
    3  = 3786
    4  = 2278
    5  = 445
    6  = 8
    7  = 17
    10 = 105
    11 = 88
    14 = 1403
    15 = 247
    16 = 2330
    17 = 165
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates4.Keys) {
    $ws4.Cells.Item($row, 6).Value = $updates4[$row]
}
